$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text-safe values (percentages, URLs, names, multi-dot numbers)
$safeValues = @{
    'D2' = '22.105.39'
    'E2' = '  +7.55%  '
    'D3' = '1.584.40'
    'E3' = '  +7.44%  '
    'E4' = '  +0.56%  '
    'E5' = '  +4.17%  '
    'E6' = '  +7.27%  '
    'E7' = '  -0.55%  '
    'E8' = '  +8.67%  '
    'E9' = '  +4.28%  '
    'E10' = '  +3.73%  '
    'E11' = '  +3.73%  '
    'E12' = '  +0.70%  '
    'E13' = '  +4.80%  '
    'E14' = '  +5.88%  '
    'E15' = '  +5.25%  '
    'E16' = '  +4.29%  '
    'D17' = '1.591.40'
    'E17' = '  +7.76%  '
    'E18' = '  +2.43%  '
    'E19' = '  +10.68%  '
    'E20' = '  +9.31%  '
    'E21' = '  +7.25%  '
    'E22' = '  +7.73%  '
    'E23' = '  +3.66%  '
    'D24' = '22.192.94'
    'E24' = '  +7.90%  '
    'E25' = '  +5.97%  '
    'E26' = '  +15.47%  '
    'E27' = '  +2.62%  '
    'E28' = '  +10.94%  '
    'D29' = '1.764.12'
    'E29' = '  +7.71%  '
    'E30' = '  +6.12%  '
    'E31' = '  +2.13%  '
    'E32' = '  +16.71%  '
    'E33' = '  +12.42%  '
    'E34' = '  +1.48%  '
    'E35' = '  +6.52%  '
    'E36' = '  +11.28%  '
    'E37' = '  +7.00%  '
    'E38' = '  +0.49%  '
    'E39' = '  +11.51%  '
    'E40' = '  +2.69%  '
    'E41' = '  +5.34%  '
    'E42' = '  +4.76%  '
    'E43' = '  +3.99%  '
    'E44' = '  +9.07%  '
    'E45' = '  +7.44%  '
    'E46' = '  +4.90%  '
    'B47' = 'Quant'
    'C47' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'E47' = '  +5.28%  '
    'B48' = 'Decentraland'
    'C48' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'E48' = '  +6.18%  '
    'E49' = '  +3.63%  '
    'E50' = '  +3.63%  '
    'E51' = '  +6.95%  '
}

# Numeric-looking values that must be forced to remain text
$textForcedValues = @{
    'D4' = '1.013'
    'D5' = '0.9966'
    'D6' = '298.21'
    'D7' = '0.3600'
    'D8' = '0.3336'
    'D9' = '41.36'
    'D10' = '1.106'
    'D11' = '0.06908'
    'D13' = '5.802'
    'D14' = '19.20'
    'D15' = '6.521'
    'D16' = '0.9976'
    'D18' = '0.00001052'
    'D19' = '0.06562'
    'D20' = '75.64'
    'D21' = '5.906'
    'D22' = '15.65'
    'D23' = '11.59'
    'D25' = '2.395'
    'D26' = '2.459'
    'D27' = '147.30'
    'D28' = '19.07'
    'D30' = '120.82'
    'D31' = '4.003'
    'D32' = '5.826'
    'D33' = '0.9094'
    'D34' = '0.08117'
    'D35' = '1.614'
    'D36' = '11.56'
    'D37' = '5.074'
    'D38' = '1.224'
    'D39' = '8.291'
    'D40' = '0.05942'
    'D41' = '0.02168'
    'D42' = '0.1969'
    'D43' = '0.9957'
    'D44' = '0.5767'
    'D45' = '3.785'
    'D46' = '12.82'
    'D47' = '124.65'
    'D48' = '0.5527'
    'D49' = '1.921'
    'D50' = '0.06716'
    'D51' = '72.06'
}

foreach ($cell in $safeValues.Keys) {
    $ws.Range($cell).Value = $safeValues[$cell]
}

foreach ($cell in $textForcedValues.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $textForcedValues[$cell]
    $rng.Style = "Normal"
}
